$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 21:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 870321
$ws.Range("C4").Value = 21604
$ws.Range("E4").Value = 736279
$ws.Range("G4").Value = 1558
$ws.Range("H4").Value = 49217

# Row 7 - Francia
$ws.Range("B7").Value = 158183
$ws.Range("C7").Value = 2239
$ws.Range("E7").Value = 94239

# Row 8 - Alemania
$ws.Range("B8").Value = 151784
$ws.Range("C8").Value = 1136
$ws.Range("E8").Value = 43080
$ws.Range("G8").Value = 89
$ws.Range("H8").Value = 5404

# Row 18 - Suiza
$ws.Range("D18").Value = 20600
$ws.Range("E18").Value = 6347

# Row 111 - Reunion
$ws.Range("B111").Value = 412
$ws.Range("C111").Value = 2
$ws.Range("E111").Value = 174

# Row 147 - Monaco
$ws.Range("D147").Value = 35
$ws.Range("E147").Value = 55
$ws.Range("G147").Value = 1
$ws.Range("H147").Value = 4

# Row 165 - Nepal
$ws.Range("B165").Value = 48
$ws.Range("C165").Value = 3
$ws.Range("E165").Value = 39
